$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"23.63000000000025"
$ws.Range("H2").Value = [double]"4.724353296277262e-16"
$ws.Range("K2").Value = [double]"41.38295220950818"
$ws.Range("L2").Value = "[33.545192989422596, 49.220711429593756]"
$ws.Range("O2").Value = [double]"1.352237078121733"
$ws.Range("P2").Value = "[1.1509738850989635, 1.553500271144502]"
$ws.Range("S2").Value = [double]"62.1733219735759"
$ws.Range("T2").Value = "[57.433395181514584, 66.9132487656372]"
$ws.Range("W2").Value = [double]"18.54446446446466"
$ws.Range("X2").Value = [double]"17.78754754754774"
$ws.Range("Y2").Value = [double]"19.30138138138159"

# Row 3 updates
$ws.Range("B3").Value = [double]"1"
$ws.Range("E3").Value = [double]"24.08000000000033"
$ws.Range("G3").Value = [double]"6.52455867111712e-12"
$ws.Range("H3").Value = [double]"2.031521064517763e-11"
$ws.Range("K3").Value = [double]"38.06233581228073"
$ws.Range("L3").Value = "[24.615236505906402, 51.509435118655055]"
$ws.Range("M3").Value = [double]"8.105795923363246e-08"
$ws.Range("N3").Value = [double]"8.105795923363246e-08"
$ws.Range("O3").Value = [double]"2.471763589310888"
$ws.Range("P3").Value = "[2.132131951084964, 2.8113952275368126]"
$ws.Range("Q3").Value = [double]"0"
$ws.Range("R3").Value = [double]"0"
$ws.Range("S3").Value = [double]"57.9170970037102"
$ws.Range("T3").Value = "[51.021475653205414, 64.81271835421498]"
$ws.Range("W3").Value = [double]"14.60708708708728"
$ws.Range("X3").Value = [double]"13.30546546546564"
$ws.Range("Y3").Value = [double]"15.90870870870893"
